$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 34 mirrors the layout of row 28 (same D/F/G/H remarks for this
# site-type combination) but records a different Site Id (column B) with
# a KK/L2100/FAIL row for site HNGRA1_BID_P47.

# Bring over row 28's formatting (fonts/borders/alignment) first.
$ws.Range("A28:H28").Copy()
$ws.Range("A34:H34").PasteSpecial(-4122)

# Column B keeps the default/no explicit style, like the source row.
$ws.Range("B34").ClearFormats()

# Now populate the values (Value2 avoids locale/formula reinterpretation).
$ws.Range("A34").Value2 = $ws.Range("A28").Value2
$ws.Range("B34").Value2 = "HNGRA1_BID_P47"
$ws.Range("C34").Value2 = $ws.Range("C28").Value2
$ws.Range("D34").Value2 = $ws.Range("D28").Value2
$ws.Range("E34").Value2 = $ws.Range("E28").Value2
$ws.Range("F34").Value2 = $ws.Range("F28").Value2
$ws.Range("G34").Value2 = $ws.Range("G28").Value2
$ws.Range("H34").Value2 = $ws.Range("H28").Value2

# Match row 28's (taller, wrapped-text) row height.
$ws.Rows.Item(34).RowHeight = 142.5

# Leave the selection/view where Excel would after typing the new last row.
$ws.Range("B34").Select()
$ws.Application.ActiveWindow.ScrollRow = 34
